# Change rat's basis components labels to alphabetical order
#
# On the "rat" worksheet, the basis-component letters T, U, V, W, X used to be
# assigned out of order (V, W, X, then T, U). This relabels them so they run
# alphabetically (T, U, V, then W, X), and updates every "Combination" formula
# cell that referenced the old letters so it uses the new letters instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rat")

# --- Runway type / Combination table (rows 3-9) ---
# Cells that used to mention the old letter V now mention the new letter T,
# and similarly for the other relabeled letters.
$ws.Range("B3").Value = "C,F,G,H,I,J,K,L,T"
$ws.Range("B4").Value = "C,F,G,H,I,J,K,L,T"
$ws.Range("B5").Value = "C,F,H,I,J,K,L,T"
$ws.Range("B7").Value = "C,Q,R,S,T"
$ws.Range("B9").Value = "C,U,V,T"
$ws.Range("B8").Value = "W,X,T"

# --- Basis parts label table (rows 30-35) ---
# Relabel the letters themselves into alphabetical order.
$ws.Range("A30").Value = "T"
$ws.Range("A31").Value = "U"
$ws.Range("A32").Value = "V"
$ws.Range("A34").Value = "W"
$ws.Range("A35").Value = "X"

# Restore the active selection on the "rat" sheet.
$ws.Activate()
$ws.Range("B5").Select()
